$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'compression pants basketball men',
    'basketball compression pants with pads',
    'basketball black',
    'sports compression pants men',
    'knee protector men',
    'athletic compression pants men',
    'basketball tights for boys',
    'protective knee pads for basketball',
    'padded knee pads basketball youth',
    'youth boys athletic leggings',
    'baseball sliding shorts youth',
    'gel knee pads for wrestling',
    'running tights capri',
    'bjj knee pads',
    'football pants with pads adult',
    'youth volleyball pads',
    'mens compression tights black',
    'little boy leggings for sports',
    'baseball men pants',
    'boys football pants with pads',
    'black knee pads volleyball',
    'knee compression for running',
    'sweat pads men',
    'basketball gear for men',
    'baseball sliding shorts',
    'black wrestling tights',
    'mens basketballs',
    'small black basketball',
    'compression baseball',
    'knee pads volleyball youth',
    'baseball pants youth black',
    'boys lacrosse pads',
    'football tights for men',
    'basketball youth',
    'working knee pads',
    'men compression clothes',
    'male compression',
    'kneepad volleyball',
    'floor hockey pads',
    'tactical pants men with knee pads',
    'sports pants for men tall',
    'lacrosse pads youth',
    'mens capri shorts',
    'sports leggings for boys',
    'knee pads for washing floors',
    'leggings under basketball shorts',
    'basketball gym',
    'capri pants mens',
    'cycling tights for men',
    'patella compression',
    'hockey knee pad',
    'mens gym pant',
    'floor pads for gym',
    'volleyball spandex',
    'nike knee pads',
    'mens volleyball knee pads',
    'mcdavid knee pads',
    'knee pad sleeve',
    'cushy knee pads',
    'men compression pants 3 4',
    'men compression pants under armour',
    'men compression pants nike',
    'mens compression tights under armour',
    'womens compression pants',
    'men compression pants xxxl',
    'men compression tights nike',
    'mens compression tights short',
    'mens compression tights thermal',
    'mens compression pants white',
    'mens compression pants navy blue',
    'mens compression pants grey',
    'mens compression pants nike',
    'mens compression pants cold gear',
    'eastbay mens compression tights',
    'womens compression pant',
    'ladies compression leggings',
    'female compression pants',
    'mcdavid basketball knee pads',
    'robo knee pads',
    'basketball long knee pads',
    'basketball mcdavid knee pads',
    'basketball nike knee pads',
    'basketball tights for kids',
    'ucla basketball youth',
    'nba basketball youth',
    'basketball 3 4 compression pants',
    'basketball youth sleeve',
    'basketball youth socks',
    'knee pad nike',
    'knee pad leggings',
    'bicycle knee pads',
    'rollerblade knee pads',
    'black leggings pants',
    'mtb knee pads',
    'kids knee pads',
    'knee pad set',
    'compression shorts with knee pads',
    'capri tights for women',
    'knee pad biking',
    'knee pads blue'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

